# Fill in previously-empty rows 5 and 6 of the utilisation report fixture
# with two new facility rows ("Crumpet" / "Scone") used to exercise the
# numerical rounding fix. Row 6 already carried the row's number formats
# (it was a blank placeholder row); row 5 needs those same formats, which
# we pick up by copying them down from row 6 before writing the values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bring row 5 up to the same cell formatting as row 6 (columns A:J only -
# row 6 additionally has a formatted-but-empty K cell that row 5 does not).
$ws.Range("A6:J6").Copy() | Out-Null
$ws.Range("A5:J5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Column A first for both new rows, then column C for both new rows -
# mirrors the order the new shared strings were actually authored in.
$ws.Range("A5").Value = "Crumpet GEF"
$ws.Range("A6").Value = "Scone GEF"
$ws.Range("C5").Value = "Crumpet exporter"
$ws.Range("C6").Value = "Scone exporter"

# Row 5 - Crumpet GEF facility
$ws.Range("B5").Value = 20001371
$ws.Range("D5").Value = "GBP"
$ws.Range("E5").Value = 7000000
$ws.Range("F5").Value = 3938753.8
$ws.Range("G5").Value = 777
$ws.Range("H5").Value = 456
$ws.Range("I5").Value = "GBP"
$ws.Range("J5").Value = "GBP"

# Row 6 - Scone GEF facility
$ws.Range("B6").Value = 20001371
$ws.Range("D6").Value = "GBP"
$ws.Range("E6").Value = 770000
$ws.Range("F6").Value = 761579.37
$ws.Range("G6").Value = 777
$ws.Range("H6").Value = 456.77
$ws.Range("I6").Value = "GBP"
$ws.Range("J6").Value = "GBP"
